# 3.1.1 support save enum to json file
# Update the st_level sheet: the generator config JSON (C1) now declares
# IsGenEnum=true with a Path, the old "Resources/subFolder" cell (D1) goes
# away, and two sample enum-with-value rows are recorded (L3 / M3).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("st_level")

# Record example enum-with-value strings first so the shared-string table
# gets them appended in this order before the JSON config string below.
$ws.Range("L3").Value = "ENUM1:5"
$ws.Range("M3").Value = "ENUM3:67"

# Update the generator config JSON in C1 to turn on enum generation with a
# path, and drop the now-unused "Resources/subFolder" cell in D1.
$ws.Range("C1").Value = '{"IsStringId":false,"IsGenItemClass":true,"JSONName":"st_levelJSON","IsGenEnum":true,"Path":"toanstt"}'
$ws.Range("D1").ClearContents()

# Widen column L (enum-with-value samples need more room) and move the
# active selection to G3, matching the saved view state.
$ws.Range("L1").ColumnWidth = 10.166666666666666
$ws.Range("G3").Select()
